{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Rewrites the \"pain points\" and \"solution/goals\" bullet lists in the\n// Executive Summary from the project-management scenario to the movie-\n// ticket-booking-website scenario, and drops the last bullet of each list.\n\nconst replacements = [\n  { oldText: \"The pain points/problems/needs/happiness (V\u1ea5n \u0111\u1ec1 ,kh\u00f3 kh\u0103n khi th\u1ef1c hi\u1ec7n d\u1ef1 \u00e1n) :\", newText: \"The pain points/problems/needs/happiness (V\u1ea5n \u0111\u1ec1) :\" },\n  { oldText: \"          B\u1ea5t \u0111\u1ed3ng \u00fd ki\u1ebfn v\u1edbi \u0111\u1ed1i t\u00e1c\", newText: \"          B\u1ea5t l\u1ee3i trong vi\u1ec7c ph\u1ea3i ph\u1ea3i t\u1edbi r\u1ea1p chi\u1ebfu phim \u0111\u1ec3 mua v\u00e9 , xem l\u1ecbch chi\u1ebfu phim\" },\n  { oldText: \"          Thi\u1ebfu nh\u00e2n l\u1ef1c => kh\u00f4ng ho\u00e0n th\u00e0nh k\u1ecbp h\u1ea1n c\u1ee7a d\u1ef1 \u00e1n => \u0111\u1ec1n b\u00f9 h\u1ee3p \u0111\u1ed3ng\", newText: \"          T\u1ed1n th\u1eddi gian x\u1ebfp h\u00e0ng \u0111\u1ec3 mua v\u00e9 \" },\n  { oldText: \"          Ph\u00e2n chia c\u00f4ng vi\u1ec7c ch\u01b0a r\u00f5 r\u00e0ng\", newText: \"          Sau khi mua v\u00e9 xong th\u00ec ph\u1ea3i ch\u1edd t\u1edbi xu\u1ea5t chi\u1ebfu\" },\n  { oldText: \"          Trong qu\u00e1 tr\u00ecnh th\u1ef1c hi\u1ec7n d\u1ef1 \u00e1n nh\u00e0 \u0111\u1ea7u t\u01b0 giao y\u00eau c\u1ea7u m\u1edbi v\u00e0 g\u1ea5p r\u00fat\", newText: \"        Kh\u00e1c h\u00e0ng t\u1ed1n th\u1eddi gian khi ph\u1ea3i xem danh s\u00e1ch xu\u1ea5t chi\u1ebfu (b\u00ecnh th\u01b0\u1eddng danh s\u00e1ch xu\u1ea5t chi\u1ebfu hi\u1ec3n th\u1ecb theo slide )\" },\n  { oldText: \"         Ph\u00e1t sinh th\u00eam chi ph\u00ed\", newText: null },\n  { oldText: \"          Ki\u1ec3m so\u00e1t ch\u1eb7t kh\u00e2u nh\u00e2n s\u1ef1 , b\u00e1o c\u00e1o t\u00ecnh h\u00ecnh ti\u1ebfn tr\u00ecnh d\u1ef1 \u00e1n m\u1ed7i tu\u1ea7n \", newText: \"          T\u1ea1o web \u0111\u1eb7t v\u00e9 xem phim ,ng\u01b0\u1eddi d\u00f9ng c\u00f3 th\u1ec3 \u0111\u1eb7t v\u00e9 tr\u1ef1c ti\u1ebfp m\u00e0 kh\u00f4ng c\u1ea7n ph\u1ea3i \u0111i t\u1edbi r\u1ea1p  v\u00e0 x\u1ebfp h\u00e0ng\" },\n  { oldText: \"          C\u00f3 nh\u00e2n l\u1ef1c d\u1ef1 tr\u00f9 \u0111\u1ec3 d\u1ec5 d\u00e0ng h\u1ed7 tr\u1ee3 ngu\u1ed3n nh\u00e2n l\u1ef1c b\u1ecb m\u1ea5t \u0111\u1ec3 \u0111i \u0111\u00fang ti\u1ebfn \u0111\u1ed9\", newText: \"      Ng\u01b0\u1eddi d\u00f9ng c\u00f3 th\u1ec3 xem c\u00e1c xu\u1ea5t chi\u1ebfu m\u00e0 kh\u00f4ng c\u1ea7n t\u1edbi r\u1ea1p chi\u1ebfu phim\" },\n  { oldText: \"          Team BA c\u00f3 k\u1ef9 n\u0103ng t\u1ed1t \u0111\u1ec3 \u0111\u00e0m ph\u00e1n v\u1edbi nh\u00e0 \u0111\u1ea7u t\u01b0 \u0111\u1ec3 \u0111\u1ea9y d\u1ef1 \u00e1n \u0111i v\u00e0o quy tr\u00ecnh b\u1eaft \u0111\u1ea7u \", newText: \" Ch\u1ec9 c\u1ea7n \u0111\u0103t v\u00e9 xem th\u00f4ng tin th\u1eddi gian v\u00e0 t\u1edbi \u0111\u00fang gi\u1edd l\u00e0 c\u00f3 th\u1ec3 xem phim li\u1ec1n ,ng\u01b0\u1eddi d\u00f9ng kh\u00f4ng c\u1ea7n ph\u1ea3i ch\u1edd \" },\n  { oldText: \"           Ng\u01b0\u1eddi qu\u1ea3n l\u00fd n\u1eafm b\u1eaft r\u00f5 r\u00e0ng y\u00eau c\u1ea7u m\u1edbi c\u1ee7a nh\u00e0 \u0111\u1ea7u t\u01b0 v\u00e0 n\u1ebfu g\u1ea5p r\u00fat th\u00ec th\u00eam nh\u00e2n l\u1ef1c d\u1ef1 tr\u00f9 \u0111\u1ec3 h\u1ed7 tr\u1ee3 cho team  \", newText: \"Ng\u01b0\u1eddi d\u00f9ng th\u00edch th\u1ec3 lo\u1ea1i phim n\u00e0o th\u00ec c\u00f3 th\u1ec3 d\u1ec5 d\u00e0ng search th\u1ec3 lo\u1ea1i phim \" },\n  { oldText: \"            Ki\u1ec3m so\u00e1t chi ph\u00ed n\u1ebfu ph\u00e1t sinh ph\u1ea3i bi\u1ebft ph\u00e1t sinh t\u1eeb \u0111\u00e2u v\u00e0 \u0111\u01b0a gi\u1ea3i ph\u00e1p x\u1eed l\u00fd vi ph\u1ea1m \u0111\u00f3\", newText: null },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Paragraphs slated for deletion (processed after all text edits so the\n// index-independent, proxy-based deletes don't disturb the other edits).\nconst toDelete = [];\n\nfor (const para of paragraphs.items) {\n  const match = replacements.find((r) => para.text === r.oldText);\n  if (!match) continue;\n  if (match.newText === null) {\n    toDelete.push(para);\n  } else {\n    para.insertText(match.newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $d resolve against the open document ($word.ActiveDocument).\n#\n# Rewrites the \"pain points\" and \"solution/goals\" bullet lists in the\n# Executive Summary from the project-management scenario to the movie-\n# ticket-booking-website scenario, and drops the last bullet of each list.\n\n$d = $word.ActiveDocument\n\n# --- 1) Text replacements (paragraph count stays unchanged for these) -----\n$replacements = @(\n  @{ Old = \"The pain points/problems/needs/happiness (V\u1ea5n \u0111\u1ec1 ,kh\u00f3 kh\u0103n khi th\u1ef1c hi\u1ec7n d\u1ef1 \u00e1n) :\"; New = \"The pain points/problems/needs/happiness (V\u1ea5n \u0111\u1ec1) :\" },\n  @{ Old = \"          B\u1ea5t \u0111\u1ed3ng \u00fd ki\u1ebfn v\u1edbi \u0111\u1ed1i t\u00e1c\"; New = \"          B\u1ea5t l\u1ee3i trong vi\u1ec7c ph\u1ea3i ph\u1ea3i t\u1edbi r\u1ea1p chi\u1ebfu phim \u0111\u1ec3 mua v\u00e9 , xem l\u1ecbch chi\u1ebfu phim\" },\n  @{ Old = \"          Thi\u1ebfu nh\u00e2n l\u1ef1c => kh\u00f4ng ho\u00e0n th\u00e0nh k\u1ecbp h\u1ea1n c\u1ee7a d\u1ef1 \u00e1n => \u0111\u1ec1n b\u00f9 h\u1ee3p \u0111\u1ed3ng\"; New = \"          T\u1ed1n th\u1eddi gian x\u1ebfp h\u00e0ng \u0111\u1ec3 mua v\u00e9 \" },\n  @{ Old = \"          Ph\u00e2n chia c\u00f4ng vi\u1ec7c ch\u01b0a r\u00f5 r\u00e0ng\"; New = \"          Sau khi mua v\u00e9 xong th\u00ec ph\u1ea3i ch\u1edd t\u1edbi xu\u1ea5t chi\u1ebfu\" },\n  @{ Old = \"          Trong qu\u00e1 tr\u00ecnh th\u1ef1c hi\u1ec7n d\u1ef1 \u00e1n nh\u00e0 \u0111\u1ea7u t\u01b0 giao y\u00eau c\u1ea7u m\u1edbi v\u00e0 g\u1ea5p r\u00fat\"; New = \"        Kh\u00e1c h\u00e0ng t\u1ed1n th\u1eddi gian khi ph\u1ea3i xem danh s\u00e1ch xu\u1ea5t chi\u1ebfu (b\u00ecnh th\u01b0\u1eddng danh s\u00e1ch xu\u1ea5t chi\u1ebfu hi\u1ec3n th\u1ecb theo slide )\" },\n  @{ Old = \"          Ki\u1ec3m so\u00e1t ch\u1eb7t kh\u00e2u nh\u00e2n s\u1ef1 , b\u00e1o c\u00e1o t\u00ecnh h\u00ecnh ti\u1ebfn tr\u00ecnh d\u1ef1 \u00e1n m\u1ed7i tu\u1ea7n \"; New = \"          T\u1ea1o web \u0111\u1eb7t v\u00e9 xem phim ,ng\u01b0\u1eddi d\u00f9ng c\u00f3 th\u1ec3 \u0111\u1eb7t v\u00e9 tr\u1ef1c ti\u1ebfp m\u00e0 kh\u00f4ng c\u1ea7n ph\u1ea3i \u0111i t\u1edbi r\u1ea1p  v\u00e0 x\u1ebfp h\u00e0ng\" },\n  @{ Old = \"          C\u00f3 nh\u00e2n l\u1ef1c d\u1ef1 tr\u00f9 \u0111\u1ec3 d\u1ec5 d\u00e0ng h\u1ed7 tr\u1ee3 ngu\u1ed3n nh\u00e2n l\u1ef1c b\u1ecb m\u1ea5t \u0111\u1ec3 \u0111i \u0111\u00fang ti\u1ebfn \u0111\u1ed9\"; New = \"      Ng\u01b0\u1eddi d\u00f9ng c\u00f3 th\u1ec3 xem c\u00e1c xu\u1ea5t chi\u1ebfu m\u00e0 kh\u00f4ng c\u1ea7n t\u1edbi r\u1ea1p chi\u1ebfu phim\" },\n  @{ Old = \"          Team BA c\u00f3 k\u1ef9 n\u0103ng t\u1ed1t \u0111\u1ec3 \u0111\u00e0m ph\u00e1n v\u1edbi nh\u00e0 \u0111\u1ea7u t\u01b0 \u0111\u1ec3 \u0111\u1ea9y d\u1ef1 \u00e1n \u0111i v\u00e0o quy tr\u00ecnh b\u1eaft \u0111\u1ea7u \"; New = \" Ch\u1ec9 c\u1ea7n \u0111\u0103t v\u00e9 xem th\u00f4ng tin th\u1eddi gian v\u00e0 t\u1edbi \u0111\u00fang gi\u1edd l\u00e0 c\u00f3 th\u1ec3 xem phim li\u1ec1n ,ng\u01b0\u1eddi d\u00f9ng kh\u00f4ng c\u1ea7n ph\u1ea3i ch\u1edd \" },\n  @{ Old = \"           Ng\u01b0\u1eddi qu\u1ea3n l\u00fd n\u1eafm b\u1eaft r\u00f5 r\u00e0ng y\u00eau c\u1ea7u m\u1edbi c\u1ee7a nh\u00e0 \u0111\u1ea7u t\u01b0 v\u00e0 n\u1ebfu g\u1ea5p r\u00fat th\u00ec th\u00eam nh\u00e2n l\u1ef1c d\u1ef1 tr\u00f9 \u0111\u1ec3 h\u1ed7 tr\u1ee3 cho team  \"; New = \"Ng\u01b0\u1eddi d\u00f9ng th\u00edch th\u1ec3 lo\u1ea1i phim n\u00e0o th\u00ec c\u00f3 th\u1ec3 d\u1ec5 d\u00e0ng search th\u1ec3 lo\u1ea1i phim \" }\n)\n\nforeach ($rep in $replacements) {\n  $rng = $d.Content\n  $rng.Find.Execute($rep.Old, $false, $false, $false, $false, $false, $true, 1, $false, $rep.New, 2)\n}\n\n# --- 2) Whole-paragraph deletions ------------------------------------------\n# Walk paragraphs back-to-front so deleting one doesn't shift the index of\n# paragraphs still queued for inspection/removal.\n$toRemove = @(\n  \"         Ph\u00e1t sinh th\u00eam chi ph\u00ed\",\n  \"            Ki\u1ec3m so\u00e1t chi ph\u00ed n\u1ebfu ph\u00e1t sinh ph\u1ea3i bi\u1ebft ph\u00e1t sinh t\u1eeb \u0111\u00e2u v\u00e0 \u0111\u01b0a gi\u1ea3i ph\u00e1p x\u1eed l\u00fd vi ph\u1ea1m \u0111\u00f3\"\n)\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n  $para = $d.Paragraphs($i)\n  $text = $para.Range.Text\n  if ($text.EndsWith(\"`r\")) { $text = $text.Substring(0, $text.Length - 1) }\n  if ($toRemove -contains $text) {\n    $para.Range.Delete()\n  }\n}\n"}
